$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CityResaleNum")

$row = 20

# Columns A and D look like a date / a plain number to Excel's auto-detection,
# so force them to be stored as text (matching the rest of the column),
# then drop the temporary "@" number format so no stray style sticks to the cell.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2023-06-04"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = "12:17:01"
$ws.Cells.Item($row, 3).Value = "Sunday"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "22"
$ws.Cells.Item($row, 4).ClearFormats()

$ws.Cells.Item($row, 5).Value = 120937
$ws.Cells.Item($row, 6).Value = 134117
$ws.Cells.Item($row, 7).Value = 158977
$ws.Cells.Item($row, 8).Value = 130051
$ws.Cells.Item($row, 9).Value = 174708
$ws.Cells.Item($row, 10).Value = 112914
$ws.Cells.Item($row, 11).Value = 199903
$ws.Cells.Item($row, 12).Value = 218799
$ws.Cells.Item($row, 13).Value = 172143
$ws.Cells.Item($row, 14).Value = 119086
$ws.Cells.Item($row, 15).Value = 38210
$ws.Cells.Item($row, 16).Value = 34764
$ws.Cells.Item($row, 17).Value = 50254
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 36665
$ws.Cells.Item($row, 20).Value = -1
